$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price (D) and Hora (G) columns we will touch as Text,
# so Excel keeps the numeric-looking strings as text (matching the
# existing cell type) instead of silently converting them to numbers.
$ws.Range("D2:D7").NumberFormat = "@"
$ws.Range("D9:D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D40:D45").NumberFormat = "@"
$ws.Range("D47:D49").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Coin / Link / Price / Volume(1h) updates
$ws.Range("D2").Value = '245.88'
$ws.Range("D3").Value = '22.11'
$ws.Range("D4").Value = '5.361'
$ws.Range("D5").Value = '0.05934'
$ws.Range("D6").Value = '3.393'
$ws.Range("D7").Value = '6.397'
$ws.Range("D9").Value = '0.9681'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1432'
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '0.03512'
$ws.Range("E11").Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("D12").Value = '0.07406'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.03044'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.09412'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").Value = '4.001'
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = '0.001590'
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").Value = '0.04825'
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").Value = '0.0005912'
$ws.Range("E18").Value = '17OneONE'
$ws.Range("D19").Value = '0.005976'
$ws.Range("D20").Value = '0.004087'
$ws.Range("D21").Value = '0.0009879'
$ws.Range("D23").Value = '3.740'
$ws.Range("D40").Value = '0.03947'
$ws.Range("D41").Value = '0.006511'
$ws.Range("D42").Value = '0.1075'
$ws.Range("D43").Value = '0.002700'
$ws.Range("D44").Value = '0.005480'
$ws.Range("D45").Value = '0.00005295'
$ws.Range("D47").Value = '0.6602'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOINBestin24h'
$ws.Range("D48").Value = '0.04572'
$ws.Range("D49").Value = '0.00002101'

# Hora column refresh (3 -> 4)
$ws.Range("G2").Value = '4'
$ws.Range("G3").Value = '4'
$ws.Range("G4").Value = '4'
$ws.Range("G5").Value = '4'
$ws.Range("G6").Value = '4'
$ws.Range("G7").Value = '4'
$ws.Range("G8").Value = '4'
$ws.Range("G9").Value = '4'
$ws.Range("G10").Value = '4'
$ws.Range("G11").Value = '4'
$ws.Range("G12").Value = '4'
$ws.Range("G13").Value = '4'
$ws.Range("G14").Value = '4'
$ws.Range("G15").Value = '4'
$ws.Range("G16").Value = '4'
$ws.Range("G17").Value = '4'
$ws.Range("G18").Value = '4'
$ws.Range("G19").Value = '4'
$ws.Range("G20").Value = '4'
$ws.Range("G21").Value = '4'
$ws.Range("G22").Value = '4'
$ws.Range("G23").Value = '4'
$ws.Range("G24").Value = '4'
$ws.Range("G25").Value = '4'
$ws.Range("G26").Value = '4'
$ws.Range("G27").Value = '4'
$ws.Range("G28").Value = '4'
$ws.Range("G29").Value = '4'
$ws.Range("G30").Value = '4'
$ws.Range("G31").Value = '4'
$ws.Range("G32").Value = '4'
$ws.Range("G33").Value = '4'
$ws.Range("G34").Value = '4'
$ws.Range("G35").Value = '4'
$ws.Range("G36").Value = '4'
$ws.Range("G37").Value = '4'
$ws.Range("G38").Value = '4'
$ws.Range("G39").Value = '4'
$ws.Range("G40").Value = '4'
$ws.Range("G41").Value = '4'
$ws.Range("G42").Value = '4'
$ws.Range("G43").Value = '4'
$ws.Range("G44").Value = '4'
$ws.Range("G45").Value = '4'
$ws.Range("G46").Value = '4'
$ws.Range("G47").Value = '4'
$ws.Range("G48").Value = '4'
$ws.Range("G49").Value = '4'
$ws.Range("G50").Value = '4'
$ws.Range("G51").Value = '4'
Write-Output "Updated symbol list"
